$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.051.28"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "'1.645.74"
$ws.Range("E3").Value = "  +0.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.53%  "

# Row 5
$ws.Range("D5").Value = "'216.02"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("E6").Value = "  +0.43%  "

# Row 7
$ws.Range("E7").Value = "  +0.54%  "

# Row 8
$ws.Range("E8").Value = "  +0.73%  "

# Row 9
$ws.Range("E9").Value = "  +0.53%  "

# Row 10
$ws.Range("D10").Value = "'19.60"
$ws.Range("E10").Value = "  +0.11%  "

# Row 11
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.27"
$ws.Range("E12").Value = "  +0.72%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.660.33"
$ws.Range("E13").Value = "  +0.65%  "

# Row 14
$ws.Range("E14").Value = "  +0.30%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0₃0764"
$ws.Range("E15").Value = "  +0.91%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'63.53"
$ws.Range("E16").Value = "  +1.68%  "

# Row 17
$ws.Range("D17").Value = "'26.064.39"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("E18").Value = "  +0.50%  "

# Row 19
$ws.Range("D19").Value = "'194.55"
$ws.Range("E19").Value = "  +0.66%  "

# Row 20
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("E21").Value = "  +0.32%  "

# Row 22
$ws.Range("E22").Value = "  -0.67%  "

# Row 23
$ws.Range("D23").Value = "'0.133"
$ws.Range("E23").Value = "  +4.97%  "

# Row 24
$ws.Range("E24").Value = "  -0.23%  "

# Row 25
$ws.Range("E25").Value = "  +0.44%  "

# Row 26
$ws.Range("D26").Value = "'143.70"
$ws.Range("E26").Value = "  -0.20%  "

# Row 27
$ws.Range("D27").Value = "'6.89"
$ws.Range("E27").Value = "  +0.64%  "

# Row 28
$ws.Range("D28").Value = "'15.52"
$ws.Range("E28").Value = "  +0.41%  "

# Row 29
$ws.Range("D29").Value = "'1.25"
$ws.Range("E29").Value = "  +0.58%  "

# Row 30
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").Value = "'3.27"
$ws.Range("E31").Value = "  +1.64%  "

# Row 32
$ws.Range("E32").Value = "  -0.11%  "

# Row 33
$ws.Range("D33").Value = "'1.55"
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("E34").Value = "  +1.44%  "

# Row 35
$ws.Range("E35").Value = "  +0.47%  "

# Row 36
$ws.Range("D36").Value = "'1.132.14"
$ws.Range("E36").Value = "  -0.38%  "

# Row 37
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("E38").Value = "  +0.44%  "

# Row 39
$ws.Range("E39").Value = "  +0.40%  "

# Row 40
$ws.Range("D40").Value = "'5.46"
$ws.Range("E40").Value = "  +0.76%  "

# Row 41
$ws.Range("D41").Value = "'99.02"
$ws.Range("E41").Value = "  -0.24%  "

# Row 42
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("E43").Value = "  +1.72%  "

# Row 44
$ws.Range("D44").Value = "'56.59"
$ws.Range("E44").Value = "  +0.26%  "

# Row 45
$ws.Range("D45").Value = "'1.50"
$ws.Range("E45").Value = "  +3.15%  "

# Row 46
$ws.Range("E46").Value = "  -1.18%  "

# Row 47
$ws.Range("D47").Value = "'7.81"
$ws.Range("E47").Value = "  +2.07%  "

# Row 48
$ws.Range("E48").Value = "  +0.00%  "

# Row 49
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("D50").Value = "'0.0952"
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("E51").Value = "  +2.94%  "

Write-Host "All updates applied."
